$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '42.945.03'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.40%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.553.33'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.03%  '

$ws.Range("E4").Value = '  -0.16%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '304.78'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.69%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '98.27'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +5.71%  '

$ws.Range("E7").Value = '  +0.11%  '

$ws.Range("E8").Value = '  +0.01%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.546'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.84%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '37.05'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.49%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0821'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.61%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.76'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.16%  '

$ws.Range("E13").Value = '  +5.62%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.944.01'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.17%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.500.57'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.37%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '15.01'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +6.26%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.878'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.38%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '42.961.96'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.47%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.73'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +4.55%  '

$ws.Range("E20").Value = '  +1.48%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.60'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.42%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '71.98'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.08%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '254.24'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.00%  '

$ws.Range("E24").Value = '  +0.60%  '

$ws.Range("E25").Value = '  -2.39%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '28.05'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -4.99%  '

$ws.Range("E27").Value = '  -0.13%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.19'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.30%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '38.06'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.91%  '

$ws.Range("E30").Value = '  -1.00%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.15'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.74%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '158.59'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.74%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.17'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.77%  '

$ws.Range("E34").Value = '  -0.62%  '

$ws.Range("B35").Value = 'Celestia'
$ws.Range("C35").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '19.37'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +13.85%  '

$ws.Range("B36").Value = 'Hedera'
$ws.Range("C36").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0805'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.11%  '

$ws.Range("E37").Value = '  -2.05%  '

$ws.Range("B38").Value = 'Kaspa'
$ws.Range("C38").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.116'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.33%  '

$ws.Range("B39").Value = 'EnergySwap'
$ws.Range("C39").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '25.79'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +10.02%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.120'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.26%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.10'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +33.10%  '

$ws.Range("B42").Value = 'NEARProtocol'
$ws.Range("C42").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.43'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.83%  '

$ws.Range("B43").Value = 'RenderToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.90'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.10%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0307'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.93%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.082.41'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.09%  '

$ws.Range("E46").Value = '  +0.05%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '86.53'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.03%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.00'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.27%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.801.49'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.14%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '74.89'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +8.19%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '103.42'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.42%  '
